# Auto-generated edit script applying market-price refresh updates
# to the Twintania_Profits workbook (per scheduled-runner diff).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1979586.8
$ws.Range("J17").Value = 2291934.8
$ws.Range("L17").Value = 6875804.399999999
$ws.Range("N17").Value = -6876140.399999999
$ws.Range("H74").Value = 5142.067
$ws.Range("I74").Value = 5176.091
$ws.Range("K74").Value = 5176.091
$ws.Range("M74").Value = -4240.091
$ws.Range("H77").Value = 5142.067
$ws.Range("I77").Value = 5176.091
$ws.Range("K77").Value = 25880.455
$ws.Range("M77").Value = -21200.455
$ws.Range("H80").Value = 436160.03
$ws.Range("J80").Value = 835304.9399999999
$ws.Range("L80").Value = 2505914.82
$ws.Range("N80").Value = -2507910.82
$ws.Range("H82").Value = 2299.6667
$ws.Range("I82").Value = 2299.6667
$ws.Range("K82").Value = 6899.000100000001
$ws.Range("M82").Value = -6493.000100000001
$ws.Range("H83").Value = 436160.03
$ws.Range("J83").Value = 835304.9399999999
$ws.Range("L83").Value = 7517744.459999999
$ws.Range("N83").Value = -7527728.459999999
$ws.Range("H85").Value = 2299.6667
$ws.Range("I85").Value = 2299.6667
$ws.Range("K85").Value = 6899.000100000001
$ws.Range("M85").Value = -5495.000100000001
$ws.Range("H88").Value = 1807.8823
$ws.Range("I88").Value = 1742.875
$ws.Range("K88").Value = 1742.875
$ws.Range("M88").Value = -1336.875
$ws.Range("H91").Value = 1807.8823
$ws.Range("I91").Value = 1742.875
$ws.Range("K91").Value = 1742.875
$ws.Range("M91").Value = -338.875
$ws.Range("H100").Value = 47259.086
$ws.Range("I100").Value = 62175.53
$ws.Range("J100").Value = 4995.8335
$ws.Range("K100").Value = 62175.53
$ws.Range("L100").Value = 4995.8335
$ws.Range("M100").Value = -61634.53
$ws.Range("N100").Value = -6077.8335
$ws.Range("H112").Value = 2251.6
$ws.Range("J112").Value = 2275.4814
$ws.Range("L112").Value = 6826.4442
$ws.Range("N112").Value = -9042.4442

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 3808
$ws.Range("I74").Value = 3182.2273
$ws.Range("K74").Value = 3182.2273
$ws.Range("M74").Value = -2308.2273
$ws.Range("H77").Value = 3808
$ws.Range("I77").Value = 3182.2273
$ws.Range("K77").Value = 15911.1365
$ws.Range("M77").Value = -11543.1365

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H9").Value = 34949
$ws.Range("J9").Value = 34949
$ws.Range("L9").Value = 34949
$ws.Range("N9").Value = -35285

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2587.0208
$ws.Range("I31").Value = 1503.28
$ws.Range("J31").Value = 3765
$ws.Range("K31").Value = 1503.28
$ws.Range("L31").Value = 3765
$ws.Range("M31").Value = -1208.28
$ws.Range("N31").Value = -4355
$ws.Range("H34").Value = 2587.0208
$ws.Range("I34").Value = 1503.28
$ws.Range("J34").Value = 3765
$ws.Range("K34").Value = 1503.28
$ws.Range("L34").Value = 3765
$ws.Range("M34").Value = -1301.28
$ws.Range("N34").Value = -4169
$ws.Range("H58").Value = 3459.4468
$ws.Range("I58").Value = 1441.5385
$ws.Range("J58").Value = 5957.8096
$ws.Range("K58").Value = 1441.5385
$ws.Range("L58").Value = 5957.8096
$ws.Range("M58").Value = -1238.5385
$ws.Range("N58").Value = -6363.8096
$ws.Range("H136").Value = 3459.4468
$ws.Range("I136").Value = 1441.5385
$ws.Range("J136").Value = 5957.8096
$ws.Range("K136").Value = 4324.6155
$ws.Range("L136").Value = 17873.4288
$ws.Range("M136").Value = -1774.6155
$ws.Range("N136").Value = -22973.4288

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1575.4054
$ws.Range("J5").Value = 1689.36
$ws.Range("L5").Value = 5068.08
$ws.Range("N5").Value = -5292.08
$ws.Range("H23").Value = 178.72728
$ws.Range("I23").Value = 88
$ws.Range("J23").Value = 198.88889
$ws.Range("K23").Value = 264
$ws.Range("L23").Value = 596.6666700000001
$ws.Range("M23").Value = -29
$ws.Range("N23").Value = -1066.66667
$ws.Range("H110").Value = 14749
$ws.Range("I110").Value = 14749
$ws.Range("K110").Value = 44247
$ws.Range("M110").Value = -40157
$ws.Range("H131").Value = 1892.8448
$ws.Range("I131").Value = 734.63635
$ws.Range("J131").Value = 2163.9148
$ws.Range("K131").Value = 2203.90905
$ws.Range("L131").Value = 6491.7444
$ws.Range("M131").Value = 2836.09095
$ws.Range("N131").Value = -16571.7444
$ws.Range("H135").Value = 1575.4054
$ws.Range("J135").Value = 1689.36
$ws.Range("L135").Value = 15204.24
$ws.Range("N135").Value = -20274.24

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3414.5715
$ws.Range("J80").Value = 3367
$ws.Range("L80").Value = 3367
$ws.Range("N80").Value = -5363
$ws.Range("H83").Value = 3414.5715
$ws.Range("J83").Value = 3367
$ws.Range("L83").Value = 16835
$ws.Range("N83").Value = -26819
$ws.Range("H102").Value = 1605.2222
$ws.Range("I102").Value = 1359.5625
$ws.Range("K102").Value = 1359.5625
$ws.Range("M102").Value = 262.4375
$ws.Range("H126").Value = 4499.75
$ws.Range("I126").Value = 4595.364
$ws.Range("J126").Value = 3448
$ws.Range("K126").Value = 13786.092
$ws.Range("L126").Value = 10344
$ws.Range("M126").Value = -11316.092
$ws.Range("N126").Value = -15284

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 2404.9
$ws.Range("I61").Value = 1800
$ws.Range("K61").Value = 1800
$ws.Range("M61").Value = -1598
$ws.Range("H113").Value = 2404.9
$ws.Range("I113").Value = 1800
$ws.Range("K113").Value = 1800
$ws.Range("M113").Value = 370
$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()
$ws.Range("H124").Value = 0
$ws.Range("J124").Value = 0
$ws.Range("L124").Value = 0
$ws.Range("N124").ClearContents()
$ws.Range("H132").Value = 5459.2856
$ws.Range("I132").Value = 3681.25
$ws.Range("K132").Value = 11043.75
$ws.Range("M132").Value = -8513.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 22500
$ws.Range("J41").Value = 15000
$ws.Range("L41").Value = 15000
$ws.Range("N41").Value = -15780
$ws.Range("H132").Value = 23555.844
$ws.Range("I132").Value = 14221.37
$ws.Range("J132").Value = 37557.555
$ws.Range("K132").Value = 42664.11
$ws.Range("L132").Value = 112672.665
$ws.Range("M132").Value = -40134.11
$ws.Range("N132").Value = -117732.665
